$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Insert a new row at position 11 ("selectRepeatevery" / 1), pushing the
# "selectweekdaysfriday" row (and everything after it) down by one.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "selectRepeatevery"
$ws.Range("B11").Value = 1

# The "weekly" value cell no longer carries the highlighted style.
$ws.Range("B10").Style = "Normal"

# Old row 11 ("selectweekdaysfriday"/"click") is now row 12; its answer
# becomes "FRI" and also drops the highlighted style.
$ws.Range("B12").Value = "FRI"
$ws.Range("B12").Style = "Normal"

# Make the Input sheet the active / selected tab, with D18 selected,
# mirroring the NewLoanInput sheet losing that state.
[void]$ws.Range("D18").Select()
$ws.Activate()
